$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the existing bold/bordered style from column A (row 45) down to the
# new rows 46-61 so the new index cells match the style used by existing rows.
$ws.Range("A45").Copy() | Out-Null
$ws.Range("A46:A61").PasteSpecial(-4122) | Out-Null

# Update cell values (rows 2-45 changed in place; rows 46-61 are newly added).
$ws.Cells.Item(2, 2).Value = 3641.900387983458
$ws.Cells.Item(2, 3).Value = 3
$ws.Cells.Item(2, 4).Value = 0.02192257005422071
$ws.Cells.Item(2, 5).Value = 0.002291987910894773
$ws.Cells.Item(3, 2).Value = 3591.307804308955
$ws.Cells.Item(3, 3).Value = 5
$ws.Cells.Item(3, 4).Value = 0.02663869637456516
$ws.Cells.Item(3, 5).Value = 0.004091081574291559
$ws.Cells.Item(4, 2).Value = 3590.077663255123
$ws.Cells.Item(4, 4).Value = 0.02366025429992787
$ws.Cells.Item(4, 5).Value = 0.00521842977910457
$ws.Cells.Item(5, 2).Value = 3550.048164038397
$ws.Cells.Item(5, 4).Value = 0.01444569267815698
$ws.Cells.Item(5, 5).Value = 0.004887875377967198
$ws.Cells.Item(6, 2).Value = 3502.960474122052
$ws.Cells.Item(6, 4).Value = 0.02772589099459388
$ws.Cells.Item(6, 5).Value = 0.005105999225098543
$ws.Cells.Item(7, 2).Value = 3431.102829766671
$ws.Cells.Item(7, 4).Value = 0.03086408861912291
$ws.Cells.Item(7, 5).Value = 0.004395878790919296
$ws.Cells.Item(8, 2).Value = 3368.080569266842
$ws.Cells.Item(8, 4).Value = 0.02309599615133376
$ws.Cells.Item(8, 5).Value = 0.004986305948373994
$ws.Cells.Item(9, 2).Value = 3307.944159957081
$ws.Cells.Item(9, 4).Value = 0.0308079837024549
$ws.Cells.Item(9, 5).Value = 0.004920212833460962
$ws.Cells.Item(10, 2).Value = 3248.901609482813
$ws.Cells.Item(10, 4).Value = 0.01478001044755087
$ws.Cells.Item(10, 5).Value = 0.002939826292938453
$ws.Cells.Item(11, 2).Value = 3194.870240440814
$ws.Cells.Item(11, 4).Value = 0.02020765854304886
$ws.Cells.Item(11, 5).Value = 0.003248773928114032
$ws.Cells.Item(12, 2).Value = 3055.02847882692
$ws.Cells.Item(12, 4).Value = 0.02128466388230035
$ws.Cells.Item(12, 5).Value = 0.002940387905617116
$ws.Cells.Item(13, 2).Value = 2975.873467839226
$ws.Cells.Item(13, 4).Value = 0.02770523131760127
$ws.Cells.Item(13, 5).Value = 0.003396105129249617
$ws.Cells.Item(14, 2).Value = 2914.997110421229
$ws.Cells.Item(14, 4).Value = 0.01792353509950995
$ws.Cells.Item(14, 5).Value = 0.00307352635875652
$ws.Cells.Item(15, 2).Value = 2818.003645176748
$ws.Cells.Item(15, 4).Value = 0.031864291236701
$ws.Cells.Item(15, 5).Value = 0.005147245597063278
$ws.Cells.Item(16, 2).Value = 2754.078401999958
$ws.Cells.Item(16, 4).Value = 0.01811789826080882
$ws.Cells.Item(16, 5).Value = 0.003247970228214284
$ws.Cells.Item(17, 2).Value = 2668.081568129679
$ws.Cells.Item(17, 4).Value = 0.02456670250621414
$ws.Cells.Item(17, 5).Value = 0.003536048011891281
$ws.Cells.Item(18, 2).Value = 2602.076549741646
$ws.Cells.Item(18, 4).Value = 0.02409938030112982
$ws.Cells.Item(18, 5).Value = 0.004949263546595271
$ws.Cells.Item(19, 2).Value = 2710.045077107244
$ws.Cells.Item(19, 3).Value = 6
$ws.Cells.Item(19, 4).Value = 0.01051513105247635
$ws.Cells.Item(19, 5).Value = 0.002734785003882226
$ws.Cells.Item(20, 2).Value = 2531.033588374834
$ws.Cells.Item(20, 4).Value = 0.03772381185211433
$ws.Cells.Item(20, 5).Value = 0.004760357633397656
$ws.Cells.Item(21, 2).Value = 2420.949001992897
$ws.Cells.Item(21, 3).Value = 6
$ws.Cells.Item(21, 4).Value = 0.03446884615579243
$ws.Cells.Item(21, 5).Value = 0.003840073856848417
$ws.Cells.Item(22, 2).Value = 2368.880808504633
$ws.Cells.Item(22, 4).Value = 0.02011961836457342
$ws.Cells.Item(22, 5).Value = 0.00426744314970525
$ws.Cells.Item(23, 2).Value = 2322.153129981494
$ws.Cells.Item(23, 3).Value = 6
$ws.Cells.Item(23, 4).Value = 0.02719728923297746
$ws.Cells.Item(23, 5).Value = 0.005100844512277635
$ws.Cells.Item(24, 2).Value = 2278.055618030553
$ws.Cells.Item(24, 3).Value = 3
$ws.Cells.Item(24, 4).Value = 0.02909305941778273
$ws.Cells.Item(24, 5).Value = 0.007019101399252999
$ws.Cells.Item(25, 2).Value = 2204.358204112928
$ws.Cells.Item(25, 3).Value = 3
$ws.Cells.Item(25, 4).Value = 0.03006435580960081
$ws.Cells.Item(25, 5).Value = 0.00625422613225518
$ws.Cells.Item(26, 2).Value = 2169.050429654774
$ws.Cells.Item(26, 4).Value = 0.01895810366595726
$ws.Cells.Item(26, 5).Value = 0.004582627388520497
$ws.Cells.Item(27, 2).Value = 2132.06047007272
$ws.Cells.Item(27, 3).Value = 6
$ws.Cells.Item(27, 4).Value = 0.01845117016752726
$ws.Cells.Item(27, 5).Value = 0.004024333868752509
$ws.Cells.Item(28, 2).Value = 2076.921596470905
$ws.Cells.Item(28, 3).Value = 6
$ws.Cells.Item(28, 4).Value = 0.02963968111937093
$ws.Cells.Item(28, 5).Value = 0.004811233526620771
$ws.Cells.Item(29, 2).Value = 2022.42484984847
$ws.Cells.Item(29, 3).Value = 2
$ws.Cells.Item(29, 4).Value = 0.09979395805368245
$ws.Cells.Item(29, 5).Value = 0.02515718207778839
$ws.Cells.Item(30, 2).Value = 1970.787546092246
$ws.Cells.Item(30, 4).Value = 0.01596931668162964
$ws.Cells.Item(30, 5).Value = 0.002368410314732851
$ws.Cells.Item(31, 2).Value = 1921.155009467111
$ws.Cells.Item(31, 3).Value = 1
$ws.Cells.Item(31, 4).Value = 0.1725349402278302
$ws.Cells.Item(31, 5).Value = 0.01914584612281324
$ws.Cells.Item(32, 2).Value = 1879.848117975443
$ws.Cells.Item(32, 3).Value = 6
$ws.Cells.Item(32, 4).Value = 0.03148268152150961
$ws.Cells.Item(32, 5).Value = 0.0036565711116972
$ws.Cells.Item(33, 2).Value = 1837.108717688487
$ws.Cells.Item(33, 3).Value = 6
$ws.Cells.Item(33, 4).Value = 0.01896815130355755
$ws.Cells.Item(33, 5).Value = 0.003359282187239796
$ws.Cells.Item(34, 2).Value = 1776.304464390151
$ws.Cells.Item(34, 3).Value = 3
$ws.Cells.Item(34, 4).Value = 0.03154956860255804
$ws.Cells.Item(34, 5).Value = 0.00565554358419269
$ws.Cells.Item(35, 2).Value = 1654.530638235761
$ws.Cells.Item(35, 3).Value = 6
$ws.Cells.Item(35, 4).Value = 0.02259084688376933
$ws.Cells.Item(35, 5).Value = 0.003733829649748973
$ws.Cells.Item(36, 2).Value = 1620.253637000357
$ws.Cells.Item(36, 4).Value = 0.1858667272389751
$ws.Cells.Item(36, 5).Value = 0.02489450488727895
$ws.Cells.Item(37, 2).Value = 1571.342367480472
$ws.Cells.Item(37, 3).Value = 3
$ws.Cells.Item(37, 4).Value = 0.03429536003964168
$ws.Cells.Item(37, 5).Value = 0.003177792878461232
$ws.Cells.Item(38, 2).Value = 1533.028047880304
$ws.Cells.Item(38, 3).Value = 6
$ws.Cells.Item(38, 4).Value = 0.02065408367364211
$ws.Cells.Item(38, 5).Value = 0.002545582131048032
$ws.Cells.Item(39, 2).Value = 1481.739655643233
$ws.Cells.Item(39, 3).Value = 6
$ws.Cells.Item(39, 4).Value = 0.01453373587559134
$ws.Cells.Item(39, 5).Value = 0.003517793588354757
$ws.Cells.Item(40, 2).Value = 1374.586558670756
$ws.Cells.Item(40, 3).Value = 1
$ws.Cells.Item(40, 4).Value = 0.5935097690188837
$ws.Cells.Item(40, 5).Value = 0.06037726556651042
$ws.Cells.Item(41, 2).Value = 1322.294118779407
$ws.Cells.Item(41, 3).Value = 6
$ws.Cells.Item(41, 4).Value = 0.009717756258840842
$ws.Cells.Item(41, 5).Value = 0.003402487004729026
$ws.Cells.Item(42, 2).Value = 1300.028177532532
$ws.Cells.Item(42, 3).Value = 6
$ws.Cells.Item(42, 4).Value = 0.01312659772080421
$ws.Cells.Item(42, 5).Value = 0.003566527071155214
$ws.Cells.Item(43, 2).Value = 1256.98613235723
$ws.Cells.Item(43, 4).Value = 0.07444713341309173
$ws.Cells.Item(43, 5).Value = 0.008392975417379975
$ws.Cells.Item(44, 2).Value = 1212.786733327345
$ws.Cells.Item(44, 3).Value = 4
$ws.Cells.Item(44, 4).Value = 0.04849649450086303
$ws.Cells.Item(44, 5).Value = 0.007518295814826646
$ws.Cells.Item(45, 2).Value = 1171.224394255621
$ws.Cells.Item(45, 3).Value = 2
$ws.Cells.Item(45, 4).Value = 0.105190508708432
$ws.Cells.Item(45, 5).Value = 0.02320056215857191
$ws.Cells.Item(46, 1).Value = 44
$ws.Cells.Item(46, 2).Value = 1126.3
$ws.Cells.Item(46, 3).Value = 1
$ws.Cells.Item(46, 4).Value = 0.142341353077429
$ws.Cells.Item(46, 5).Value = 0.01615234951150906
$ws.Cells.Item(47, 1).Value = 45
$ws.Cells.Item(47, 2).Value = 1092.1
$ws.Cells.Item(47, 3).Value = 1
$ws.Cells.Item(47, 4).Value = 0.4607198742420117
$ws.Cells.Item(47, 5).Value = 0.07202972712886309
$ws.Cells.Item(48, 1).Value = 46
$ws.Cells.Item(48, 2).Value = 1062.101907177814
$ws.Cells.Item(48, 3).Value = 4
$ws.Cells.Item(48, 4).Value = 0.0786444411539772
$ws.Cells.Item(48, 5).Value = 0.01109265924966474
$ws.Cells.Item(49, 1).Value = 47
$ws.Cells.Item(49, 2).Value = 891.5075023775776
$ws.Cells.Item(49, 3).Value = 6
$ws.Cells.Item(49, 4).Value = 0.01707563542601586
$ws.Cells.Item(49, 5).Value = 0.003599946535452295
$ws.Cells.Item(50, 1).Value = 48
$ws.Cells.Item(50, 2).Value = 808.109281249695
$ws.Cells.Item(50, 3).Value = 0
$ws.Cells.Item(50, 4).Value = 5.655299060451213
$ws.Cells.Item(50, 5).Value = 0.5557969655556487
$ws.Cells.Item(51, 1).Value = 49
$ws.Cells.Item(51, 2).Value = 749.4
$ws.Cells.Item(51, 3).Value = 2
$ws.Cells.Item(51, 4).Value = 0.4409901200150917
$ws.Cells.Item(51, 5).Value = 0.09148078207759601
$ws.Cells.Item(52, 1).Value = 50
$ws.Cells.Item(52, 2).Value = 698.4294621497356
$ws.Cells.Item(52, 3).Value = 5
$ws.Cells.Item(52, 4).Value = 0.02913580019687534
$ws.Cells.Item(52, 5).Value = 0.004895474606561604
$ws.Cells.Item(53, 1).Value = 51
$ws.Cells.Item(53, 2).Value = 662.5350823811776
$ws.Cells.Item(53, 3).Value = 3
$ws.Cells.Item(53, 4).Value = 0.02505980313866853
$ws.Cells.Item(53, 5).Value = 0.003243506251967708
$ws.Cells.Item(54, 1).Value = 52
$ws.Cells.Item(54, 2).Value = 607.646340827428
$ws.Cells.Item(54, 3).Value = 5
$ws.Cells.Item(54, 4).Value = 0.05073179914482301
$ws.Cells.Item(54, 5).Value = 0.006071340057549751
$ws.Cells.Item(55, 1).Value = 53
$ws.Cells.Item(55, 2).Value = 490.2399813553229
$ws.Cells.Item(55, 3).Value = 4
$ws.Cells.Item(55, 4).Value = 0.1482266241326562
$ws.Cells.Item(55, 5).Value = 0.03932543089233735
$ws.Cells.Item(56, 1).Value = 54
$ws.Cells.Item(56, 2).Value = 472.4903753881296
$ws.Cells.Item(56, 3).Value = 2
$ws.Cells.Item(56, 4).Value = 2.495872980536268
$ws.Cells.Item(56, 5).Value = 0.170790686070911
$ws.Cells.Item(57, 1).Value = 55
$ws.Cells.Item(57, 2).Value = 417.2
$ws.Cells.Item(57, 3).Value = 6
$ws.Cells.Item(57, 4).Value = 0.02241488784086501
$ws.Cells.Item(57, 5).Value = 0.007685104402582291
$ws.Cells.Item(58, 1).Value = 56
$ws.Cells.Item(58, 2).Value = 393.1751238336715
$ws.Cells.Item(58, 3).Value = 3
$ws.Cells.Item(58, 4).Value = 0.05630835840339664
$ws.Cells.Item(58, 5).Value = 0.01553334024921286
$ws.Cells.Item(59, 1).Value = 57
$ws.Cells.Item(59, 2).Value = 360.5
$ws.Cells.Item(59, 3).Value = 2
$ws.Cells.Item(59, 4).Value = 0.8539068883101714
$ws.Cells.Item(59, 5).Value = 0.06505957244267972
$ws.Cells.Item(60, 1).Value = 58
$ws.Cells.Item(60, 2).Value = 229.1
$ws.Cells.Item(60, 3).Value = 2
$ws.Cells.Item(60, 4).Value = 0.6789971084094479
$ws.Cells.Item(60, 5).Value = 0.0814075113775491
$ws.Cells.Item(61, 1).Value = 59
$ws.Cells.Item(61, 2).Value = 181
$ws.Cells.Item(61, 3).Value = 5
$ws.Cells.Item(61, 4).Value = 1.925789853069213
$ws.Cells.Item(61, 5).Value = 0.06428176034885699
